$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("G2").Value = 15193
$ws.Range("H2").Value = 21072
$ws.Range("J2").Value = 163.59
$ws.Range("K2").Value = 222

# Update row 3
$ws.Range("G3").Value = 15193
$ws.Range("H3").Value = 21072
$ws.Range("J3").Value = 163.59
$ws.Range("K3").Value = 222

# Move the active selection to match the diff
$ws.Range("M6").Select()
